# Updates the cryptos list: Price (D) and Volume(1h) (E) columns for a number
# of rows, per the GitHub Actions crypto-price refresh commit.
# D-column values that look like plain numbers are written with a leading
# apostrophe so Excel keeps them as text (matching the source data, which
# stores these values as text strings, not numbers); the cell style is then
# reset to Normal so the quote-prefix marker doesn't leave a stray style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.755.25"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "3.756.42"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'602.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").Value = "'169.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "3.755.14"
$ws.Range("E7").Value = "  -1.13%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("E10").Value = "  +3.44%  "
$ws.Range("D11").Value = "'6.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.71%  "
$ws.Range("D12").Value = "'0.463"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").Value = "'38.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").Value = "'0.0000247"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").Value = "4.375.55"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").Value = "3.748.75"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").Value = "68.734.94"
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("D18").Value = "'7.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "'17.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").Value = "'10.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.07%  "
$ws.Range("D22").Value = "'495.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'0.731"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.66%  "
$ws.Range("D24").Value = "'85.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").Value = "'0.0000146"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").Value = "'2.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.05%  "
$ws.Range("D27").Value = "'12.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").Value = "'10.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "'2.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.60%  "
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D32").Value = "'7.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").Value = "'32.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("D34").Value = "3.897.14"
$ws.Range("D35").Value = "3.683.23"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("D42").Value = "'437.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.71%  "
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").Value = "'8.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("D48").Value = "'40.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").Value = "2.823.02"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D51").Value = "'0.0356"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.93%  "
